# The workbook had several cells whose text contained embedded line breaks
# (manual "\n" within the cell string, e.g. from wrapped/multi-line entries
# in the shared-strings table). This edit removes those line breaks,
# replacing each one with a single space, while leaving everything else
# (including the trailing "..." truncation markers already present in the
# source data) untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet1: "Products" column (J2:J30) -----------------------------------
# All 29 data rows shared the same text "المنتجات\nالبند..." -> one space.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("J2:J30").Value = "المنتجات البند..."

# --- Sheet3: a handful of individual cells ---------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Column C ("Main Activity") rows 23 and 28.
$ws3.Range("C23").Value = "الأثاث        ..."
$ws3.Range("C28").Value = "البلاستيك     ..."

# Column I ("Business Turnover") row 26.
$ws3.Range("I26").Value = "100 ريال سعودي..."

# Column J ("Annual Export Value") rows 10 and 26.
$ws3.Range("J10").Value = "5 M ريال سعودي..."
$ws3.Range("J26").Value = "1.5 M ريال سعو..."
